# Applies the cryptos.xlsx price/volume refresh described in the commit.
# Source cells are plain text (inline strings), e.g. "67.546.78" or "  -1.03%  ".
# Excel auto-detects numeric-looking strings and would silently convert them to
# real numbers (e.g. "6.50" -> 6.5) unless the cell is explicitly formatted as
# Text ("@") before the value is assigned. We restore the style afterwards so the
# cell keeps its original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.546.78"
$ws.Range("E2").Value = "  -1.03%  "
Set-TextValue "D3" "3.761.78"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "595.52"
$ws.Range("E5").Value = "  -0.78%  "
Set-TextValue "D6" "170.73"
$ws.Range("E6").Value = "  +0.88%  "
Set-TextValue "D7" "3.760.59"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue "D9" "0.525"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +1.14%  "
Set-TextValue "D11" "6.50"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("E13").Value = "  +6.29%  "
Set-TextValue "D14" "36.66"
$ws.Range("E14").Value = "  -0.95%  "
Set-TextValue "D15" "4.391.36"
$ws.Range("E15").Value = "  -1.56%  "
Set-TextValue "D16" "3.760.86"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("E17").Value = "  +0.36%  "
Set-TextValue "D18" "67.523.86"
$ws.Range("E18").Value = "  -1.05%  "
Set-TextValue "D19" "7.21"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("E21").Value = "  -5.02%  "
Set-TextValue "D22" "468.67"
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "83.92"
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D25" "0.0000147"
$ws.Range("E25").Value = "  -8.67%  "
Set-TextValue "D26" "2.22"
$ws.Range("E26").Value = "  -0.28%  "
Set-TextValue "D27" "12.17"
$ws.Range("E27").Value = "  +0.42%  "
Set-TextValue "D28" "10.37"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("E29").Value = "  -0.17%  "
Set-TextValue "D30" "2.91"
$ws.Range("E30").Value = "  -1.50%  "
Set-TextValue "D31" "3.909.63"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("E32").Value = "  -0.14%  "
Set-TextValue "D33" "30.62"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  -3.00%  "
Set-TextValue "D35" "9.14"
$ws.Range("E35").Value = "  -3.34%  "
Set-TextValue "D36" "3.725.81"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D39" "1.00"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D40" "0.138"
$ws.Range("E40").Value = "  -1.80%  "
Set-TextValue "D41" "5.85"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("E42").Value = "  -0.11%  "
Set-TextValue "D43" "0.312"
$ws.Range("E43").Value = "  -0.68%  "
Set-TextValue "D45" "8.76"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("E47").Value = "  -2.20%  "
Set-TextValue "D48" "399.91"
$ws.Range("E48").Value = "  -4.21%  "
$ws.Range("E49").Value = "  -8.42%  "
Set-TextValue "D50" "140.10"
$ws.Range("E50").Value = "  -1.16%  "
Set-TextValue "D51" "0.0354"
$ws.Range("E51").Value = "  -0.93%  "
